$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the three address values used across the device inventory.
# Block 1 (rows 2-33):  10.6.3.56 -> 192.168.1.56
# Block 2 (rows 34-65): 10.6.3.55 -> 192.168.1.55
# Block 3 (rows 66-97): 10.6.3.54 -> 192.168.1.54
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Value = "192.168.1.56"
}
for ($r = 34; $r -le 65; $r++) {
    $ws.Cells.Item($r, 3).Value = "192.168.1.55"
}
for ($r = 66; $r -le 97; $r++) {
    $ws.Cells.Item($r, 3).Value = "192.168.1.54"
}

# Column width adjustments
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(4).ColumnWidth = 6.1640625

# Scroll position
$ws.Application.ActiveWindow.ScrollRow = 22
